# ---------------------------------------------------------------------------
# Applies the "Add files via upload" commit to AAPL.xlsx:
#   * rename Sheet1 -> Main, Sheet2 -> Model
#   * on Main: drop the "AAPL" title text and the "notes:" label, and move the
#     Price/Shares/MC/Cash/Debt/EV/PE block from columns B/E/F/G to O/P/Q
#   * on Model: re-point the three formulas that referenced Sheet1!F2/F3/F5/F6
#     at Main!P2/P3/P5/P6, and tweak a couple of cosmetic view settings
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "Main"
$ws2.Name = "Model"

# ---------------------------------------------------------------------------
# Main sheet (formerly Sheet1)
# ---------------------------------------------------------------------------

# A1 used to hold the shared string "AAPL" - blank it out but keep its (bold)
# formatting.
$ws1.Range("A1").ClearContents()

# B2 held the "notes:" label and is removed entirely.
$ws1.Range("B2").Clear()

# Text labels (Price, Shares, MC, Cash, Debt, EV, PE) move E->O. Columns past
# G already default to the same (Arial, General) style those labels used, so
# no explicit formatting is required.
$ws1.Range("O2").Value = "Price"
$ws1.Range("O3").Value = "Shares"
$ws1.Range("O4").Value = "MC"
$ws1.Range("O5").Value = "Cash"
$ws1.Range("O6").Value = "Debt"
$ws1.Range("O7").Value = "EV"
$ws1.Range("O8").Value = "PE"

# "Q125" labels move G->Q (also plain default style).
$ws1.Range("Q3").Value = "Q125"
$ws1.Range("Q5").Value = "Q125"
$ws1.Range("Q6").Value = "Q125"

# Values / formulas move F->P. Set the formula/value first, *then* the number
# format, since typing a formula that references an already-formatted cell
# can otherwise make Excel copy that cell's format onto the new one.
$ws1.Range("P2").Value = 200
$ws1.Range("P2").NumberFormat = "#,##0.00"

$ws1.Range("P3").Value = 15115.823
$ws1.Range("P3").NumberFormat = "#,##0"

$ws1.Range("P4").Formula = "=P3*P2"
$ws1.Range("P4").NumberFormat = "#,##0"

$ws1.Range("P5").Formula = "=30299+23476"
$ws1.Range("P5").NumberFormat = "#,##0"

$ws1.Range("P6").Value = 83956
$ws1.Range("P6").NumberFormat = "#,##0"

$ws1.Range("P7").Formula = "=P4+P6-P5"
$ws1.Range("P7").NumberFormat = "#,##0"

# P8 keeps the default General format, same as its old F8 location.
$ws1.Range("P8").Formula = "=P2/Model!M16"

# Now remove the old B/E/F/G cells that were moved above.
$ws1.Range("E2:G8").Clear()

# Sheet view: zoom 220 -> 115, selection A2 -> M2
$ws1.Activate()
$excel.ActiveWindow.Zoom = 115
$ws1.Range("M2").Select()

# ---------------------------------------------------------------------------
# Model sheet (formerly Sheet2)
# ---------------------------------------------------------------------------

$ws2.Range("Z22").Formula = "=NPV(Z21,N36:FK36)+Main!P5-Main!P6"
$ws2.Range("Z23").Formula = "=Z22/Main!P3"
$ws2.Range("Z24").Formula = "=Z23/Main!P2-1"

# Column Z got a bit wider and lost its "best fit" flag.
$ws2.Range("Z1").ColumnWidth = 10

# Sheet view: zoom 100 -> 130, selection Z23 -> X19, scroll position tweak.
$ws2.Activate()
$excel.ActiveWindow.Zoom = 130
$ws2.Range("X19").Select()

$wb.Application.Calculate()
Write-Host "edit complete"
